# "Usda Case Wizard Completion"
# Update the credential/value table on Sheet1 (A1:E6) to reflect the
# completed USDA case wizard run:
#   - B2 (URL)      -> https://smartnsc.com/
#   - B3 (UserName) -> Amitthakur
#   - B4 (Password) -> Aamit55555000@
# Also drop the Hyperlink cell style from C4 and E4 (they keep their
# text/hyperlink relationship, just render with the default style now).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://smartnsc.com/"
$ws.Range("B3").Value = "Amitthakur"
$ws.Range("B4").Value = "Aamit55555000@"

$ws.Range("C4").Style = "Normal"
$ws.Range("E4").Style = "Normal"

# Reflect the last active selection recorded in the saved workbook.
$ws.Range("D9").Select()
